# Add team Wins/Losses/Ties record columns (AD, AE, AF) to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels in row 1 (columns AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered/top aligned, thin border)
# by copying the format from an existing header cell (A1) onto the new headers.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill every data row (2 through 64) with the team's record for the season.
$ws.Range("AD2:AD64").Value = 50
$ws.Range("AE2:AE64").Value = 112
$ws.Range("AF2:AF64").Value = 0
